# Rebuild the literature review table (rows 1-10, cols A-F) with the
# updated headers/tool rows described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe existing values + per-cell formatting in the table range so stale
# styles/strings do not linger once the grid is rewritten below.
$ws.Range("A1:F10").Clear()

# Row 1
$ws.Range("A1").Value = "Tool name"
$ws.Range("B1").Value = "Open Source"
$ws.Range("C1").Value = "Long-Term Investment"
$ws.Range("D1").Value = "Market"
$ws.Range("E1").Value = "Stochastic Inputs"
$ws.Range("F1").Value = "Country Generalisability"

# Row 2
$ws.Range("A2").Value = "SEPIA"
$ws.Range("B2").Value = "✓"
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 16
$ws.Range("B2").Font.Color = 2236962
$ws.Range("C2").Value = "x"
$ws.Range("D2").Value = "Bilateral"
$ws.Range("D2").Font.Name = "Arial"
$ws.Range("D2").Font.Size = 16
$ws.Range("D2").Font.Color = 2236962
$ws.Range("E2").Value = "Demand"

# Row 3
$ws.Range("A3").Value = "EMCAS"
$ws.Range("B3").Value = "x"
$ws.Range("C3").Value = "✓"
$ws.Range("C3").Font.Name = "Arial"
$ws.Range("C3").Font.Size = 16
$ws.Range("C3").Font.Color = 2236962
$ws.Range("D3").Value = "✓"
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D3").Font.Size = 16
$ws.Range("D3").Font.Color = 2236962
$ws.Range("E3").Value = "Outages"
$ws.Range("F3").Value = "✓"
$ws.Range("F3").Font.Name = "Arial"
$ws.Range("F3").Font.Size = 16
$ws.Range("F3").Font.Color = 2236962

# Row 4
$ws.Range("A4").Value = "NEMSIM"
$ws.Range("B4").Value = "Unknown"
$ws.Range("C4").Value = "✓"
$ws.Range("C4").Font.Name = "Arial"
$ws.Range("C4").Font.Size = 16
$ws.Range("C4").Font.Color = 2236962
$ws.Range("D4").Value = "✓"
$ws.Range("D4").Font.Name = "Arial"
$ws.Range("D4").Font.Size = 16
$ws.Range("D4").Font.Color = 2236962
$ws.Range("E4").Value = "x"
$ws.Range("F4").Value = "x"

# Row 5
$ws.Range("A5").Value = "AMES"
$ws.Range("B5").Value = "✓"
$ws.Range("B5").Font.Name = "Arial"
$ws.Range("B5").Font.Size = 16
$ws.Range("B5").Font.Color = 2236962
$ws.Range("C5").Value = "x"
$ws.Range("D5").Value = "Day-ahead"
$ws.Range("D5").Font.Name = "Arial"
$ws.Range("D5").Font.Size = 16
$ws.Range("D5").Font.Color = 2236962
$ws.Range("E5").Value = "x"
$ws.Range("F5").Value = "x"

# Row 6
$ws.Range("A6").Value = "PowerACE"
$ws.Range("B6").Value = "x"
$ws.Range("C6").Value = "✓"
$ws.Range("C6").Font.Name = "Arial"
$ws.Range("C6").Font.Size = 16
$ws.Range("C6").Font.Color = 2236962
$ws.Range("D6").Value = "✓"
$ws.Range("D6").Font.Name = "Arial"
$ws.Range("D6").Font.Size = 16
$ws.Range("D6").Font.Color = 2236962
$ws.Range("E6").Value = "Outages/Demand"
$ws.Range("E6").Font.Name = "Arial"
$ws.Range("E6").Font.Size = 16
$ws.Range("E6").Font.Color = 2236962
$ws.Range("F6").Value = "EU"
$ws.Range("F6").Font.Name = "Arial"
$ws.Range("F6").Font.Size = 16
$ws.Range("F6").Font.Color = 2236962

# Row 7
$ws.Range("A7").Value = "MACSEM"
$ws.Range("B7").Value = "Unknown"
$ws.Range("C7").Value = "x"
$ws.Range("D7").Value = "✓"
$ws.Range("D7").Font.Name = "Arial"
$ws.Range("D7").Font.Size = 16
$ws.Range("D7").Font.Color = 2236962
$ws.Range("E7").Value = "x"
$ws.Range("F7").Value = "✓"
$ws.Range("F7").Font.Name = "Arial"
$ws.Range("F7").Font.Size = 16
$ws.Range("F7").Font.Color = 2236962

# Row 8
$ws.Range("A8").Value = "GAPEX"
$ws.Range("B8").Value = "Unknown"
$ws.Range("C8").Value = "x"
$ws.Range("D8").Value = "Day-ahead"
$ws.Range("D8").Font.Name = "Arial"
$ws.Range("D8").Font.Size = 16
$ws.Range("D8").Font.Color = 2236962
$ws.Range("E8").Value = "x"
$ws.Range("F8").Value = "✓"
$ws.Range("F8").Font.Name = "Arial"
$ws.Range("F8").Font.Size = 16
$ws.Range("F8").Font.Color = 2236962

# Row 9
$ws.Range("A9").Value = "EMLab"
$ws.Range("B9").Value = "✓"
$ws.Range("B9").Font.Name = "Arial"
$ws.Range("B9").Font.Size = 16
$ws.Range("B9").Font.Color = 2236962
$ws.Range("C9").Value = "✓"
$ws.Range("C9").Font.Name = "Arial"
$ws.Range("C9").Font.Size = 16
$ws.Range("C9").Font.Color = 2236962
$ws.Range("D9").Value = "Futures"
$ws.Range("D9").Font.Name = "Arial"
$ws.Range("D9").Font.Size = 16
$ws.Range("D9").Font.Color = 2236962
$ws.Range("E9").Value = "x"
$ws.Range("F9").Value = "✓"
$ws.Range("F9").Font.Name = "Arial"
$ws.Range("F9").Font.Size = 16
$ws.Range("F9").Font.Color = 2236962

# Row 10
$ws.Range("A10").Value = "ElecSIM"
$ws.Range("B10").Value = "✓"
$ws.Range("B10").Font.Name = "Arial"
$ws.Range("B10").Font.Size = 16
$ws.Range("B10").Font.Color = 2236962
$ws.Range("C10").Value = "✓"
$ws.Range("C10").Font.Name = "Arial"
$ws.Range("C10").Font.Size = 16
$ws.Range("C10").Font.Color = 2236962
$ws.Range("D10").Value = "Futures"
$ws.Range("D10").Font.Name = "Arial"
$ws.Range("D10").Font.Size = 16
$ws.Range("D10").Font.Color = 2236962
$ws.Range("E10").Value = "✓"
$ws.Range("E10").Font.Name = "Arial"
$ws.Range("E10").Font.Size = 16
$ws.Range("E10").Font.Color = 2236962
$ws.Range("F10").Value = "✓"
$ws.Range("F10").Font.Name = "Arial"
$ws.Range("F10").Font.Size = 16
$ws.Range("F10").Font.Color = 2236962

# Rows 2-10 use the taller 20pt row height (rows 3-5,7,10 already had it).
$ws.Rows(2).RowHeight = 20
$ws.Rows(6).RowHeight = 20
$ws.Rows(8).RowHeight = 20
$ws.Rows(9).RowHeight = 20

# Page is now portrait / A4 paper size per the new layout.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Header row selected (matches saved selection sqref="A1:XFD1").
$ws.Range("A1:XFD1").Select() | Out-Null
